# Insert a new column before column F ("PnL_per_lot") to hold the new
# "volume_weighted_avg_spread_in_USD" metric, shifting the former F:K
# columns (and their data) right to G:L, matching the author's diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a whole new column at F; Excel shifts F:K -> G:L (cells,
# formatting, column widths) automatically, and grows the sheet
# dimension from K26 to L26.
$ws.Range("F1").EntireColumn.Insert()

# New column F is 35 characters wide. ColumnWidth (COM, "characters")
# round-trips through this runtime with a constant +5/6 offset when
# written to the underlying <col width="..."> attribute, so back that
# off here to land exactly on 35.
$ws.Range("F1").EntireColumn.ColumnWidth = 34.166666666666664

# Header row for each instrument's table (rows 2, 9, 14, 21) gets the
# new column header in the freshly inserted column F.
$ws.Range("F2").Value = "volume_weighted_avg_spread_in_USD"
$ws.Range("F9").Value = "volume_weighted_avg_spread_in_USD"
$ws.Range("F14").Value = "volume_weighted_avg_spread_in_USD"
$ws.Range("F21").Value = "volume_weighted_avg_spread_in_USD"

# New per-row metric values for the new column F.
$ws.Range("F4").Value = 11.2826327013138
$ws.Range("F5").Value = 10.20510549211238
$ws.Range("F6").Value = 11.42664613780947
$ws.Range("F7").Value = 8.282943930626022

$ws.Range("F11").Value = 14.83735353158003
$ws.Range("F12").Value = 14.39382169432019

$ws.Range("F16").Value = 11.9027311541665
$ws.Range("F17").Value = 13.77828987307455
$ws.Range("F18").Value = 11.54599846330231
$ws.Range("F19").Value = 60.70948392535197

$ws.Range("F23").Value = 28.2763339998281
$ws.Range("F24").Value = 27.71893333892316
$ws.Range("F25").Value = 25.30695057091755
$ws.Range("F26").Value = 27.17329847727988
